$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.373.90"
$ws.Range("E2").Value = "  +4.77%  "
$ws.Range("D3").Value = "3.090.63"
$ws.Range("E3").Value = "  +3.49%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'556.71"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "'143.82"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +10.03%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").Value = "3.082.05"
$ws.Range("E8").Value = "  +3.49%  "
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("D10").Value = "'7.04"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +17.46%  "
$ws.Range("E11").Value = "  +5.77%  "
$ws.Range("D12").Value = "'0.459"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "'0.0000226"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.24%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'35.13"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +4.00%  "
$ws.Range("D15").Value = "3.593.10"
$ws.Range("E15").Value = "  +4.17%  "
$ws.Range("D16").Value = "64.432.73"
$ws.Range("E16").Value = "  +4.75%  "
$ws.Range("D17").Value = "3.095.92"
$ws.Range("E17").Value = "  +4.12%  "
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").Value = "'6.71"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.11%  "
$ws.Range("D20").Value = "'483.20"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.01%  "
$ws.Range("D21").Value = "'13.69"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +4.90%  "
$ws.Range("D22").Value = "'0.669"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.20%  "
$ws.Range("D23").Value = "'7.54"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +8.56%  "
$ws.Range("D24").Value = "'13.19"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +10.17%  "
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'2.81"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +4.15%  "
$ws.Range("D28").Value = "'7.96"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +4.16%  "
$ws.Range("D29").Value = "'2.05"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +8.61%  "
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("D31").Value = "'26.03"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.59%  "
$ws.Range("E32").Value = "  +2.85%  "
$ws.Range("D33").Value = "'2.43"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.87%  "
$ws.Range("D34").Value = "'5.73"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +4.70%  "
$ws.Range("D35").Value = "'55.34"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("E36").Value = "  +4.58%  "
$ws.Range("D37").Value = "'461.39"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.91%  "
$ws.Range("D38").Value = "'0.0405"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.12%  "
$ws.Range("D39").Value = "'0.0822"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.51%  "
$ws.Range("D40").Value = "3.022.80"
$ws.Range("E40").Value = "  -2.56%  "
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("D42").Value = "'8.26"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +2.82%  "
$ws.Range("D43").Value = "'2.69"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +16.38%  "
$ws.Range("D44").Value = "'28.02"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +10.55%  "
$ws.Range("D45").Value = "'0.259"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.70%  "
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "'2.07"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +7.31%  "
$ws.Range("E48").Value = "  +4.54%  "
$ws.Range("D49").Value = "'117.91"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.85%  "
$ws.Range("D50").Value = "0.0₃0515"
$ws.Range("E50").Value = "  +7.70%  "
$ws.Range("D51").Value = "'2.05"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.31%  "
